$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 597, pushing existing rows 597-615 down to 598-616.
$ws.Rows.Item(597).Insert()

# Populate the newly inserted row 597 with the new weekly record.
$ws.Cells.Item(597, 1).Value = 6
$ws.Cells.Item(597, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(597, 3).Value = "Metropolitana"
$ws.Cells.Item(597, 4).Value = "2023-05-29"
$ws.Cells.Item(597, 5).Value = 13
$ws.Cells.Item(597, 6).Value = 100112043
$ws.Cells.Item(597, 7).Value = "Pepino ensalada"
$ws.Cells.Item(597, 8).Value = "Sin especificar"
$ws.Cells.Item(597, 9).Value = "Primera"
$ws.Cells.Item(597, 10).Value = 2200
$ws.Cells.Item(597, 11).Value = 8000
$ws.Cells.Item(597, 12).Value = 9000
$ws.Cells.Item(597, 13).Value = 8455
$ws.Cells.Item(597, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(597, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(597, 16).Value = 141
$ws.Cells.Item(597, 17).Value = 60
$ws.Cells.Item(597, 18).Value = "Hortaliza"
